# Applies the updated cryptocurrency price/volume figures to Sheet1 (rows 2-51).
# Values that look like plain decimal numbers (e.g. "228.18") are written with a
# leading apostrophe so Excel stores them as text (matching the source data, which
# uses text cells throughout columns D/E, e.g. "34.406.51"), then the style is reset
# to Normal so no stray quote-prefix formatting is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.406.51'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = "'228.18"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.14%  '
$ws.Range('D6').Value = "'0.602"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.58%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = "'36.26"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.48%  '
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('D12').Value = '2.068.95'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = "'11.31"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').Value = '1.813.84'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('D17').Value = '34.410.57'
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('D18').Value = "'70.09"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('D19').Value = "'245.45"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.20%  '
$ws.Range('D20').Value = '0.0₃0788'
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').Value = "'2.26"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.88%  '
$ws.Range('D25').Value = "'171.16"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = "'8.15"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.46%  '
$ws.Range('D27').Value = "'17.53"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.48%  '
$ws.Range('E28').Value = '  +3.77%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').Value = "'3.99"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D32').Value = "'3.82"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('E33').Value = '  -1.74%  '
$ws.Range('E34').Value = '  -3.07%  '
$ws.Range('D35').Value = '1.382.56'
$ws.Range('E35').Value = '  -2.74%  '
$ws.Range('D36').Value = "'0.661"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.20%  '
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('E38').Value = '  -1.97%  '
$ws.Range('E39').Value = '  -11.57%  '
$ws.Range('D40').Value = "'82.36"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.42%  '
$ws.Range('D41').Value = "'2.81"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('D42').Value = "'0.951"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('E44').Value = '  +7.15%  '
$ws.Range('D45').Value = "'13.41"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('E46').Value = '  -2.26%  '
$ws.Range('E47').Value = '  -3.59%  '
$ws.Range('D48').Value = '1.970.08'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').Value = "'103.35"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.10%  '
$ws.Range('D51').Value = '0.0₆0123'
$ws.Range('E51').Value = '  -6.04%  '
